{"js": "// Remove the \"Mode -\" paragraph entirely (including its paragraph mark),\n// matching the target edit: the \"Median -\" paragraph stays untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text.trim() === \"Mode -\") {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the \"Mode -\" paragraph entirely (including its paragraph mark),\n# matching the target edit: the \"Median -\" paragraph stays untouched.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\", \"`n\") -eq \"Mode -\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
